# The review notes originally had two separate bullet paragraphs:
#   A) "-No existe el buscador de entrenadores. ... o apellidos."
#   B) "-No existe el botón de unirse/dejar ... o d" + [_GoBack bookmark] + "ejar actividades."
#
# The target revision removes paragraph A entirely and turns paragraph B's
# two runs into a single run holding the full sentence, with the _GoBack
# bookmark sitting right at the start of the (now sole) surviving paragraph.

$d = $word.ActiveDocument

# --- locate paragraph A ("buscador de entrenadores") and paragraph B ("botón de unirse/dejar") ---
$paraAIndex = -1
$paraBIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*No existe el buscador de entrenadores*") {
        $paraAIndex = $i
    }
    if ($text -like "*de unirse/dejar un gimnasio para customer*") {
        $paraBIndex = $i
    }
}

if ($paraAIndex -lt 1 -or $paraBIndex -lt 1) {
    throw "Could not locate the expected paragraphs (A=$paraAIndex, B=$paraBIndex)"
}

$fullText = "-No existe el botón de unirse/dejar un gimnasio para customer, con lo cual tampoco se puede ver si puede unirse o dejar actividades."

# --- remove paragraph B completely (content + its paragraph mark) ---
# Deleting A instead would leave B's paragraph mark as the "surviving" one;
# here we want A's paragraph mark/properties to be the ones that remain,
# so B (the later paragraph) is the one removed outright.
$d.Paragraphs.Item($paraBIndex).Range.Delete()

# --- rewrite paragraph A's text to the final consolidated sentence ---
$paraA = $d.Paragraphs.Item($paraAIndex)
$bodyRange = $d.Range($paraA.Range.Start, $paraA.Range.End - 1)
$bodyRange.Text = $fullText

# --- re-create the _GoBack bookmark, collapsed at the very start of the paragraph ---
$paraA = $d.Paragraphs.Item($paraAIndex)
$bookmarkRange = $d.Range($paraA.Range.Start, $paraA.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
